$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: semantic type (iaest-measure -> iaest-dimension / sdmx-dimension) ---
$ws.Range("A3").Value = "iaest-dimension:nacionalidad-pais-nombre"
$ws.Range("B3").Value = "iaest-dimension:edad-grandes-grupos"
$ws.Range("D3").Value = "sdmx-dimension:refArea"
$ws.Range("H3").Value = "iaest-dimension:nacionalidad-area-nombre"
$ws.Range("I3").Value = "iaest-dimension:sexo"

# --- Row 4: medida/dim flag follows row 3 (columns that became dimensions) ---
$ws.Range("A4").Value = "dim"
$ws.Range("B4").Value = "dim"
$ws.Range("D4").Value = "dim"
$ws.Range("H4").Value = "dim"
$ws.Range("I4").Value = "dim"

# --- Row 5: rdf type (xsd:string -> skos:Concept, except D5 -> URI-Provincia) ---
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("B5").Value = "skos:Concept"
$ws.Range("D5").Value = "URI-Provincia"
$ws.Range("H5").Value = "skos:Concept"
$ws.Range("I5").Value = "skos:Concept"

# --- Row 6 (new): mapping files for the recoded dimension columns ---
$ws.Range("A6").Value = "mapping-nacionalidad-pais-nombre.xlsx"
$ws.Range("B6").Value = "mapping-edad-grandes-grupos.xlsx"
$ws.Range("H6").Value = "mapping-nacionalidad-area-nombre.xlsx"
$ws.Range("I6").Value = "mapping-sexo.xlsx"

# Match the style used by the rest of the data rows (cellXfs index 1)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("I5").Copy()
$ws.Range("I6").PasteSpecial(-4122)
